$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Papers List")

# --- H13: new Date Reviewed cell, matching existing H-column date style ---
$ws.Range("H13").Value = 43266
$ws.Range("H13").NumberFormat = $ws.Range("G13").NumberFormat

# --- Reference cells for style propagation on the new rows ---
$bFmt = $ws.Range("B22").NumberFormat
$gFmt = $ws.Range("G23").NumberFormat

# --- Row 24: Title, Year(plain number kept in date-styled cell), Source, Date Found ---
$ws.Range("A24").Value = "Predictive environmental risk assessment of chemical mixtures: A conceptual framework"
$ws.Range("B24").Value = 2012
$ws.Range("B24").NumberFormat = $bFmt
$ws.Range("F24").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G24").Value = 43266
$ws.Range("G24").NumberFormat = $gFmt

# --- Row 25 ---
$ws.Range("A25").Value = "How well can we predict the toxicity of pesticide mixtures to aquatic life?"
$ws.Range("F25").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G25").Value = 43266
$ws.Range("G25").NumberFormat = $gFmt

# --- Row 26 ---
$ws.Range("A26").Value = "Quantifying Synergy: A Systematic Review of Mixture Toxicity Studies within Environmental Toxicology"
$ws.Range("F26").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G26").Value = 43266
$ws.Range("G26").NumberFormat = $gFmt

# --- Row 27 ---
$ws.Range("A27").Value = "Rethinking our approach to multiple stressor studies in marine environments"
$ws.Range("F27").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G27").Value = 43266
$ws.Range("G27").NumberFormat = $gFmt

# --- Row 28 ---
$ws.Range("A28").Value = "Density dependence governs when population responses to multiple stressors are magnified or mitigated"
$ws.Range("F28").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G28").Value = 43266
$ws.Range("G28").NumberFormat = $gFmt

# --- Row 29 ---
$ws.Range("A29").Value = "Interactions between effects of environmental chemicals and natural stressors: A review"
$ws.Range("F29").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G29").Value = 43266
$ws.Range("G29").NumberFormat = $gFmt

# --- Row 30 ---
$ws.Range("A30").Value = "Generalized concentration addition: A method for examining mixtures containing partial agonists"
$ws.Range("F30").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G30").Value = 43266
$ws.Range("G30").NumberFormat = $gFmt

# --- Row 31 ---
$ws.Range("A31").Value = "Predicting the synergy of multiple stress effects"
$ws.Range("F31").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G31").Value = 43266
$ws.Range("G31").NumberFormat = $gFmt

# --- Row 32 ---
$ws.Range("A32").Value = "Reconceptualizing synergism and antagonism among multiple stressors "
$ws.Range("F32").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G32").Value = 43266
$ws.Range("G32").NumberFormat = $gFmt

# --- Row 33 ---
$ws.Range("A33").Value = "Dose-Response Analysis Using R"
$ws.Range("F33").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G33").Value = 43266
$ws.Range("G33").NumberFormat = $gFmt

# --- Row 34 ---
$ws.Range("A34").Value = "Contribution of organic toxicants to multiple stress in river ecosystems"
$ws.Range("F34").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G34").Value = 43266
$ws.Range("G34").NumberFormat = $gFmt

# --- Row 35 ---
$ws.Range("A35").Value = "Generalized concentration addition approach for predicting mixture toxicity"
$ws.Range("F35").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G35").Value = 43266
$ws.Range("G35").NumberFormat = $gFmt

# --- Row 36 ---
$ws.Range("A36").Value = "An improved null model for assessing the net effects of multiple stressors on communities"
$ws.Range("F36").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G36").Value = 43266
$ws.Range("G36").NumberFormat = $gFmt

# --- Row 37 ---
$ws.Range("A37").Value = "Impacts of multiple stressors on biodiversity and ecosystem functioning: the role of species co‐tolerance"
$ws.Range("F37").Value = "Advancing understanding and prediction in multiple stressor research through a mechanistic basis for null models"
$ws.Range("G37").Value = 43266
$ws.Range("G37").NumberFormat = $gFmt

# --- View: Excel scrolled/selected this range after the edit ---
$ws.Activate() | Out-Null
$ws.Range("F24:G37").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Papers List updated"
